$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2239263803680982
$ws.Range("C2").Value = 0.5184049079754601
$ws.Range("J2").Value = 0.01840490797546012
$ws.Range("P2").Value = 0.1564417177914111
$ws.Range("S2").Value = 0.08282208588957055
$ws.Range("B3").Value = 0.005882352941176471
$ws.Range("C3").Value = 0.005882352941176471
$ws.Range("J3").Value = 0.04117647058823529
$ws.Range("P3").Value = 0.7352941176470589
$ws.Range("S3").Value = 0.2117647058823529
$ws.Range("J4").Value = 0.09090909090909091
$ws.Range("P4").Value = 0.7045454545454546
$ws.Range("S4").Value = 0.2045454545454546
$ws.Range("B6").Value = 0.0670391061452514
$ws.Range("D6").Value = 0.01675977653631285
$ws.Range("F6").Value = 0.08379888268156424
$ws.Range("J6").Value = 0.1396648044692737
$ws.Range("O6").Value = 0.01675977653631285
$ws.Range("Q6").Value = 0.2011173184357542
$ws.Range("R6").Value = 0.07262569832402235
$ws.Range("S6").Value = 0.4022346368715084
$ws.Range("B7").Value = 0.09523809523809523
$ws.Range("D7").Value = 0.005291005291005291
$ws.Range("F7").Value = 0.0582010582010582
$ws.Range("J7").Value = 0.1957671957671958
$ws.Range("O7").Value = 0.03174603174603174
$ws.Range("Q7").Value = 0.1481481481481481
$ws.Range("R7").Value = 0.07936507936507936
$ws.Range("S7").Value = 0.3862433862433862
$ws.Range("B8").Value = 0.09715639810426541
$ws.Range("D8").Value = 0.02132701421800948
$ws.Range("F8").Value = 0.05450236966824645
$ws.Range("J8").Value = 0.1327014218009479
$ws.Range("O8").Value = 0.009478672985781991
$ws.Range("Q8").Value = 0.1872037914691943
$ws.Range("R8").Value = 0.0924170616113744
$ws.Range("S8").Value = 0.4052132701421801
$ws.Range("B9").Value = 0.1060606060606061
$ws.Range("D9").Value = 0.01515151515151515
$ws.Range("F9").Value = 0.0505050505050505
$ws.Range("J9").Value = 0.1414141414141414
$ws.Range("O9").Value = 0.01515151515151515
$ws.Range("Q9").Value = 0.2121212121212121
$ws.Range("R9").Value = 0.06565656565656566
$ws.Range("S9").Value = 0.3939393939393939
$ws.Range("B10").Value = 0.1178247734138973
$ws.Range("D10").Value = 0.02190332326283988
$ws.Range("F10").Value = 0.06797583081570997
$ws.Range("J10").Value = 0.1435045317220544
$ws.Range("O10").Value = 0.01283987915407855
$ws.Range("Q10").Value = 0.195619335347432
$ws.Range("R10").Value = 0.07779456193353475
$ws.Range("S10").Value = 0.3625377643504532
$ws.Range("G11").Value = 0.1297709923664122
$ws.Range("J11").Value = 0.09541984732824428
$ws.Range("K11").Value = 0.1755725190839695
$ws.Range("L11").Value = 0.5916030534351145
$ws.Range("S11").Value = 0.007633587786259542
$ws.Range("G12").Value = 0.8209876543209876
$ws.Range("J12").Value = 0.1358024691358025
$ws.Range("K12").Value = 0.006172839506172839
$ws.Range("L12").Value = 0.03703703703703703
$ws.Range("G13").Value = 0.5813953488372093
$ws.Range("J13").Value = 0.3023255813953488
$ws.Range("S13").Value = 0.1162790697674419
$ws.Range("F15").Value = 0.004854368932038835
$ws.Range("H15").Value = 0.1067961165048544
$ws.Range("I15").Value = 0.1019417475728155
$ws.Range("J15").Value = 0.3883495145631068
$ws.Range("K15").Value = 0.0825242718446602
$ws.Range("M15").Value = 0.02912621359223301
$ws.Range("N15").Value = 0.004854368932038835
$ws.Range("O15").Value = 0.06310679611650485
$ws.Range("S15").Value = 0.2184466019417476
$ws.Range("F16").Value = 0.005050505050505051
$ws.Range("H16").Value = 0.2272727272727273
$ws.Range("I16").Value = 0.0707070707070707
$ws.Range("J16").Value = 0.3787878787878788
$ws.Range("K16").Value = 0.1313131313131313
$ws.Range("M16").Value = 0.01515151515151515
$ws.Range("O16").Value = 0.06060606060606061
$ws.Range("S16").Value = 0.1111111111111111
$ws.Range("F17").Value = 0.006756756756756757
$ws.Range("H17").Value = 0.1509009009009009
$ws.Range("I17").Value = 0.1193693693693694
$ws.Range("J17").Value = 0.4527027027027027
$ws.Range("K17").Value = 0.07432432432432433
$ws.Range("M17").Value = 0.01801801801801802
$ws.Range("O17").Value = 0.06081081081081081
$ws.Range("S17").Value = 0.1171171171171171
$ws.Range("H18").Value = 0.1978021978021978
$ws.Range("I18").Value = 0.05494505494505494
$ws.Range("J18").Value = 0.4505494505494506
$ws.Range("K18").Value = 0.1098901098901099
$ws.Range("M18").Value = 0.02197802197802198
$ws.Range("N18").Value = 0.005494505494505495
$ws.Range("O18").Value = 0.05494505494505494
$ws.Range("S18").Value = 0.1043956043956044
$ws.Range("F19").Value = 0.004035512510088781
$ws.Range("H19").Value = 0.2074253430185634
$ws.Range("I19").Value = 0.08071025020177562
$ws.Range("J19").Value = 0.3970944309927361
$ws.Range("K19").Value = 0.09362389023405973
$ws.Range("M19").Value = 0.01856335754640839
$ws.Range("O19").Value = 0.06698950766747377
$ws.Range("S19").Value = 0.1315577078288943
